# Updated cryptos list with latest price/volume data.
# Values in column D are plain-text (not numeric) fields: where the
# new text would otherwise auto-parse as a number, a leading apostrophe
# forces Excel to keep it as literal text, matching the source data's format.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.643.33"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").Value = "1.724.56"
$ws.Range("E3").Value = "  -0.37%  "
$ws.Range("D4").Value = "'0.9987"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "'241.37"
$ws.Range("E5").Value = "  -1.01%  "
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("D7").Value = "'0.4924"
$ws.Range("E7").Value = "  +0.75%  "
$ws.Range("D8").Value = "'0.2608"
$ws.Range("E8").Value = "  -0.45%  "
$ws.Range("D9").Value = "'0.06214"
$ws.Range("E9").Value = "  +0.52%  "
$ws.Range("D10").Value = "1.728.95"
$ws.Range("E10").Value = "  -0.13%  "
$ws.Range("D11").Value = "'15.83"
$ws.Range("E11").Value = "  +2.12%  "
$ws.Range("D12").Value = "'0.06988"
$ws.Range("E12").Value = "  -0.50%  "
$ws.Range("D13").Value = "'0.6085"
$ws.Range("E13").Value = "  +1.15%  "
$ws.Range("D14").Value = "'4.496"
$ws.Range("E14").Value = "  -1.28%  "
$ws.Range("D15").Value = "'77.21"
$ws.Range("E15").Value = "  -0.20%  "
$ws.Range("D16").Value = "'0.9991"
$ws.Range("E16").Value = "  -0.22%  "
$ws.Range("D17").Value = "26.474.85"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("D18").Value = "'0.9989"
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("D19").Value = "'0.000007208"
$ws.Range("E19").Value = "  +1.72%  "
$ws.Range("D20").Value = "'11.38"
$ws.Range("E20").Value = "  -0.22%  "
$ws.Range("D21").Value = "1.951.25"
$ws.Range("E21").Value = "  -0.21%  "
$ws.Range("D22").Value = "'4.460"
$ws.Range("D23").Value = "'8.537"
$ws.Range("E23").Value = "  -0.73%  "
$ws.Range("D24").Value = "'5.085"
$ws.Range("E24").Value = "  -1.90%  "
$ws.Range("D25").Value = "'137.84"
$ws.Range("E25").Value = "  -0.69%  "
$ws.Range("D26").Value = "'15.34"
$ws.Range("E26").Value = "  +0.45%  "
$ws.Range("D27").Value = "'1.762"
$ws.Range("E27").Value = "  +2.72%  "
$ws.Range("D28").Value = "'1.383"
$ws.Range("E28").Value = "  -2.16%  "
$ws.Range("D29").Value = "'106.32"
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("D30").Value = "'3.910"
$ws.Range("E30").Value = "  -1.61%  "
$ws.Range("D31").Value = "'0.07964"
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("D32").Value = "'3.661"
$ws.Range("E32").Value = "  -1.08%  "
$ws.Range("D33").Value = "'0.04490"
$ws.Range("E33").Value = "  -0.62%  "
$ws.Range("D34").Value = "'0.9977"
$ws.Range("E34").Value = "  -0.31%  "
$ws.Range("D35").Value = "'2.612"
$ws.Range("E35").Value = "  -0.19%  "
$ws.Range("D36").Value = "'0.9999"
$ws.Range("E36").Value = "  -0.17%  "
$ws.Range("D37").Value = "'0.6237"
$ws.Range("E37").Value = "  -0.28%  "
$ws.Range("D38").Value = "'0.9311"
$ws.Range("E38").Value = "  +2.56%  "
$ws.Range("D39").Value = "'2.046"
$ws.Range("E39").Value = "  +2.59%  "
$ws.Range("D40").Value = "'2.417"
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("D41").Value = "'0.9996"
$ws.Range("D42").Value = "'0.01509"
$ws.Range("E42").Value = "  +1.33%  "
$ws.Range("D43").Value = "'5.548"
$ws.Range("E43").Value = "  +1.82%  "
$ws.Range("D44").Value = "'99.22"
$ws.Range("E44").Value = "  -1.29%  "
$ws.Range("D45").Value = "'0.3845"
$ws.Range("E45").Value = "  -0.59%  "
$ws.Range("D46").Value = "'6.851"
$ws.Range("E46").Value = "  +2.63%  "
$ws.Range("D47").Value = "'0.1155"
$ws.Range("E47").Value = "  -0.21%  "
$ws.Range("D48").Value = "'0.05384"
$ws.Range("E48").Value = "  +0.35%  "
$ws.Range("D49").Value = "'7.824"
$ws.Range("E49").Value = "  +1.49%  "
$ws.Range("D50").Value = "'30.23"
$ws.Range("E50").Value = "  -0.34%  "
$ws.Range("D51").Value = "'51.56"
$ws.Range("E51").Value = "  +1.02%  "
